# Update odds values on Sheet1 to reflect the latest FlashScore data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 updates ---
$ws.Range("G3").Value = 1.57
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 6.5
$ws.Range("N3").Value = 7.5
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("AU3").Value = 10
$ws.Range("AW3").Value = 7.5

# --- Row 5 updates ---
$ws.Range("G5").Value = 1.83
$ws.Range("H5").Value = 3.35
$ws.Range("I5").Value = 4.1
$ws.Range("J5").Value = 2.37
$ws.Range("L5").Value = 4.4
$ws.Range("P5").Value = 2.72
$ws.Range("Q5").Value = 2
$ws.Range("S5").Value = 1.39
$ws.Range("T5").Value = 2.57
$ws.Range("U5").Value = 1.87
$ws.Range("Y5").Value = 8.5
$ws.Range("Z5").Value = 15
$ws.Range("AA5").Value = 16
$ws.Range("AB5").Value = 32
$ws.Range("AH5").Value = 10.5
$ws.Range("AI5").Value = 22
$ws.Range("AK5").Value = 65
$ws.Range("AL5").Value = 40
$ws.Range("AO5").Value = 9
$ws.Range("AP5").Value = 18
$ws.Range("AQ5").Value = 32
$ws.Range("AR5").Value = 65
$ws.Range("AT5").Value = 2.52
$ws.Range("AU5").Value = 7.3
$ws.Range("AW5").Value = 5.8
$ws.Range("AX5").Value = 23
$ws.Range("AZ5").Value = 120
$ws.Range("BA5").Value = 150
